$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 2; $row -le 88; $row++) {
    $cell = $ws.Cells.Item($row, 4)  # Column D = Speaker
    $val = $cell.Value2
    if ($val -eq "MOLLY MCNINCH") {
        $cell.Value = "T"
    } elseif ($val -eq "STUDENT") {
        $cell.Value = "S"
    }
}
